$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "Henrique Batalha-Filho"
$ws.Range("B11").Value = "henrique.batalha@outlook.com"
$ws.Range("A12").Value = "Júlia de Lima Carvalho"
$ws.Range("B12").Value = "julialima.carvalho18@gmail.com"
$ws.Range("A13").Value = "Emilio de Lanna Neto"
$ws.Range("B13").Value = "emiliolanna@gmail.com"
$ws.Range("A18").Value = "Diego José Santana Silva"
$ws.Range("B18").Value = "jose.santana@ufms.br"
$ws.Range("A19").Value = "Larissa Bortoli de Souza"
$ws.Range("B19").Value = "bortoli.larissa16@gmail.com"
$ws.Range("A55").Value = "Fernanda Freitas de Oliveira"
$ws.Range("B55").Value = "fernandaoliveira@ufpr.br"
$ws.Range("A56").Value = "Fabricius Domingos"
$ws.Range("B56").Value = "fabriciusmaia@gmail.com"
$ws.Range("A57").Value = "Talita Helen Bombardelli Gomig"
$ws.Range("B57").Value = "talitahbg@ufpr.br"
$ws.Range("A58").Value = "Rodrigo Barbosa Gonçalves"
$ws.Range("B58").Value = "goncalvesrb@gmail.com"
$ws.Range("A59").Value = "Juliano Morimoto"
$ws.Range("B59").Value = "juliano.morimoto@abdn.ac.uk"
$ws.Range("A60").Value = "Ricardo Lehtonen Rodrigues Souza"
$ws.Range("B60").Value = "lehtonen@ufpr.br"
$ws.Range("A61").Value = "Mayara Pereira Neves"
$ws.Range("B61").Value = "mayara-nevesbio@hotmail.com"
$ws.Range("A62").Value = "Matheus M. A. Salles"
$ws.Range("B62").Value = "matheusm.salles@hotmail.com"
$ws.Range("A63").Value = "Paulo Henrique Mueller"
$ws.Range("B63").Value = "profmueller@gmail.com"
$ws.Range("A65").Value = "Thieres Tayroni Martins da Silva"
$ws.Range("B65").Value = "thierestayroni@gmail.com"
$ws.Range("A66").Value = "Fabrício Rodrigues dos Santos"
$ws.Range("B66").Value = "fsantos.ufmg@gmail.com"
$ws.Range("A67").Value = "Lucas Bleicher"
$ws.Range("B67").Value = "lbleicher@gmail.com"
$ws.Range("A68").Value = "Isaac Rafael Freitas Borges"
$ws.Range("B68").Value = "isaacborges966@gmail.com"
$ws.Range("A69").Value = "ANA CECILIA HOLLER DEL PRETTE"
$ws.Range("B69").Value = "anacecilia.holler@gmail.com"
$ws.Range("A70").Value = "Gabriel Costa Santos"
$ws.Range("B70").Value = "gabrielscosta90@gmail.com"
$ws.Range("A71").Value = "Rafael Félix de Magalhães"
$ws.Range("B71").Value = "rafaelmagalhaes@ufsj.edu.br"
$ws.Range("A72").Value = "Ramon Moreira Fernandes"
$ws.Range("B72").Value = "ramonmf360@gmail.com"
$ws.Range("A73").Value = "Welignton Clarindo"
$ws.Range("B73").Value = "well.clarindo@ufv.br"
$ws.Range("A74").Value = "Alessandro Marques De Oliveira"
$ws.Range("B74").Value = "biolessandro@gmail.com"
$ws.Range("A75").Value = "Mariana Fonseca Rossi"
$ws.Range("B75").Value = "mfonsecarossi@gmail.com"
$ws.Range("A76").Value = "Henrique Caldeira Costa"
$ws.Range("B76").Value = "ccostah@gmail.com"
$ws.Range("A77").Value = "André Yves"
$ws.Range("B77").Value = "andreyves7@gmail.com"
$ws.Range("A78").Value = "Comissão Avaliadora"
$ws.Range("B78").Value = "angela_portella@hotmail.com"
$ws.Range("A79").Value = "RAFAEL FILGUEIRA JORGE"
$ws.Range("B79").Value = "rafajorgebio@gmail.com"
$ws.Range("A80").Value = "Ana Carolina Martins Junqueira"
$ws.Range("B80").Value = "anacmj@gmail.com"
$ws.Range("A81").Value = "Jose Ricardo Miras Mermudes"
$ws.Range("B81").Value = "jrmermudes@gmail.com"
$ws.Range("A82").Value = "Carlos Eduardo Guerra Schrago"
$ws.Range("B82").Value = "guerra@biologia.ufrj.br"
$ws.Range("A83").Value = "Miguel Godinho Alvares"
$ws.Range("B83").Value = "miguelgodinhoalvares@gmail.com"
$ws.Range("A84").Value = "Dener Soares Da Costa Junior"
$ws.Range("B84").Value = "denerdacosta12@gmail.com"
$ws.Range("A85").Value = "Bruno Loreto de Aragão Pedroso"
$ws.Range("B85").Value = "bruno.loreto.aragao@hotmail.com"
$ws.Range("A86").Value = "Claudia Augusta de Moraes Russo"
$ws.Range("B86").Value = "claurusso@hotmail.com"
$ws.Range("A87").Value = "Beatriz Mello Carvalho"
$ws.Range("B87").Value = "biaumello@gmail.com"
$ws.Range("A88").Value = "Gabriela Ferreira Mota"
$ws.Range("B88").Value = "motafgabriela@gmail.com"
$ws.Range("A89").Value = "Luiza Silva Anselmini"
$ws.Range("B89").Value = "miniansel.lu@gmail.com"
$ws.Range("A90").Value = "Lucas Pereira da Rocha"
$ws.Range("B90").Value = "lucasrocha700@gmail.com"
$ws.Range("A91").Value = "Melissa Bars Closel"
$ws.Range("B91").Value = "melissabars@gmail.com"
$ws.Range("A92").Value = "Breno Michelon Seixas"
$ws.Range("B92").Value = "breno.mseixas@usp.br"
$ws.Range("A93").Value = "Leonardo Maurici Borges"
$ws.Range("B93").Value = "aquitemcaqui@gmail.com"
$ws.Range("A94").Value = "Rafael Fernandes Barduzzi"
$ws.Range("B94").Value = "rfbarduzzi@gmail.com"
$ws.Range("A95").Value = "João Pedro Fujita"
$ws.Range("B95").Value = "joaopedrofujita@estudante.ufscar.br"
$ws.Range("A96").Value = "Lina Maria Ameida Silva"
$ws.Range("B96").Value = "linamas@gmail.com"
$ws.Range("A97").Value = "Paulo Aecyo Francisco da Silva"
$ws.Range("B97").Value = "pauloaecyo_1997@hotmail.com"
$ws.Range("A98").Value = "Pedro Paulo Goulart Taucci"
$ws.Range("B98").Value = "pedrotaucce@gmail.com"
$ws.Range("A99").Value = "Lucas Albuquerque dos Santos"
$ws.Range("B99").Value = "lucasabqsto@gmail.com"
$ws.Range("A100").Value = "João Claudio de Sousa Nascimento"
$ws.Range("B100").Value = "j.claudionasci@gmail.com"
$ws.Range("A101").Value = "Nicolle Souza Leto"
$ws.Range("B101").Value = "letonicolle@gmail.com"
$ws.Range("A102").Value = "Leonardo Duarte Santos"
$ws.Range("B102").Value = "santosldbio@gmail.com"
$ws.Range("A103").Value = "Vera Nisaka Solferini"
$ws.Range("B103").Value = "veras@unicamp.br"
$ws.Range("A104").Value = "Julia Nader Acquaviva"
$ws.Range("B104").Value = "julianader95@gmail.com"
$ws.Range("A105").Value = "Ana Claudia Lessinger"
$ws.Range("B105").Value = "lessinger@ufscar.br"
$ws.Range("A106").Value = "Marcelo Duarte"
$ws.Range("B106").Value = "mduartes@usp.br"
$ws.Range("A107").Value = "Weverton dos Santos Azevedo"
$ws.Range("B107").Value = "weverton.azevedo@hotmail.com"
$ws.Range("A108").Value = "Rafaela Velloso Missagia"
$ws.Range("B108").Value = "rafaelamissagia@gmail.com"
$ws.Range("A109").Value = "Diogo Melo"
$ws.Range("B109").Value = "diogro@gmail.com"
$ws.Range("A110").Value = "Ana Paula Assis"
$ws.Range("B110").Value = "paulaassis@ib.usp.br"
$ws.Range("A111").Value = "katarine nogueira norbertino"
$ws.Range("B111").Value = "katarinenn@outlook.com"
$ws.Range("A112").Value = "Pedro L. Godoy"
$ws.Range("B112").Value = "pedrolorenagodoy@gmail.com"
$ws.Range("A113").Value = "Gabriela Procópio Camacho"
$ws.Range("B113").Value = "gpcamacho@usp.br"
$ws.Range("A114").Value = "Mariana Mira Vasconcellos"
$ws.Range("B114").Value = "marimiravasc@gmail.com"
$ws.Range("A115").Value = "Ivan Sergio Nunes Silva Filho"
$ws.Range("B115").Value = "ivan.nunes@unesp.br"
$ws.Range("A116").Value = "Thais Helena Condez"
$ws.Range("B116").Value = "thacondez@gmail.com"
$ws.Range("A117").Value = "Jônatas Gomes Santos"
$ws.Range("B117").Value = "gomesjonatas21@gmail.com"
$ws.Range("A118").Value = "Gabriela Procopio Leite"
$ws.Range("B118").Value = "gabrielaprocopio3@gmail.com"
